# Data refresh for cryptos.xlsx: updates Price (D) and Volume(1h) (E) values
# for all coin rows, and shifts several coin name/link/price/volume rows
# (rows 27-51) to reflect the new ranking snapshot, matching the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "64.762.71"
$ws.Cells.Item(2, 5).Value = "  -3.59%  "
$ws.Cells.Item(3, 4).Value = "3.339.05"
$ws.Cells.Item(3, 5).Value = "  -4.38%  "
$ws.Cells.Item(4, 5).Value = "  -0.14%  "
$ws.Cells.Item(5, 4).Value = "'182.29"
$ws.Cells.Item(5, 5).Value = "  -8.52%  "
$ws.Cells.Item(6, 4).Value = "'534.89"
$ws.Cells.Item(6, 5).Value = "  -2.75%  "
$ws.Cells.Item(7, 5).Value = "  +0.81%  "
$ws.Cells.Item(8, 4).Value = "3.335.47"
$ws.Cells.Item(8, 5).Value = "  -4.29%  "
$ws.Cells.Item(9, 5).Value = "  +0.03%  "
$ws.Cells.Item(10, 4).Value = "'0.618"
$ws.Cells.Item(10, 5).Value = "  -5.01%  "
$ws.Cells.Item(11, 4).Value = "'58.83"
$ws.Cells.Item(11, 5).Value = "  -6.76%  "
$ws.Cells.Item(12, 4).Value = "'0.135"
$ws.Cells.Item(12, 5).Value = "  -5.49%  "
$ws.Cells.Item(13, 4).Value = "'0.0000262"
$ws.Cells.Item(14, 5).Value = "  -5.87%  "
$ws.Cells.Item(15, 4).Value = "3.872.58"
$ws.Cells.Item(15, 5).Value = "  -4.64%  "
$ws.Cells.Item(16, 4).Value = "3.341.95"
$ws.Cells.Item(16, 5).Value = "  -4.48%  "
$ws.Cells.Item(17, 5).Value = "  -4.44%  "
$ws.Cells.Item(18, 4).Value = "64.690.30"
$ws.Cells.Item(18, 5).Value = "  -3.44%  "
$ws.Cells.Item(19, 4).Value = "'17.64"
$ws.Cells.Item(19, 5).Value = "  -3.43%  "
$ws.Cells.Item(20, 4).Value = "'11.23"
$ws.Cells.Item(20, 5).Value = "  -4.22%  "
$ws.Cells.Item(21, 4).Value = "'0.968"
$ws.Cells.Item(21, 5).Value = "  -4.90%  "
$ws.Cells.Item(22, 4).Value = "'378.21"
$ws.Cells.Item(22, 5).Value = "  -2.71%  "
$ws.Cells.Item(23, 4).Value = "'3.84"
$ws.Cells.Item(23, 5).Value = "  -3.32%  "
$ws.Cells.Item(24, 4).Value = "'11.28"
$ws.Cells.Item(24, 5).Value = "  -7.01%  "
$ws.Cells.Item(25, 4).Value = "'81.30"
$ws.Cells.Item(25, 5).Value = "  -1.35%  "
$ws.Cells.Item(26, 5).Value = "  +3.16%  "
$ws.Cells.Item(27, 2).Value = "ImmutableX"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(27, 4).Value = "'2.71"
$ws.Cells.Item(27, 5).Value = "  -3.09%  "
$ws.Cells.Item(28, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(28, 4).Value = "'11.57"
$ws.Cells.Item(28, 5).Value = "  -4.68%  "
$ws.Cells.Item(29, 2).Value = "Filecoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(29, 4).Value = "'8.47"
$ws.Cells.Item(29, 5).Value = "  -3.39%  "
$ws.Cells.Item(30, 2).Value = "EthereumClassic"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(30, 4).Value = "'29.26"
$ws.Cells.Item(30, 5).Value = "  -5.22%  "
$ws.Cells.Item(31, 2).Value = "Bittensor"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(31, 4).Value = "'657.82"
$ws.Cells.Item(31, 5).Value = "  -3.22%  "
$ws.Cells.Item(32, 2).Value = "NEARProtocol"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(32, 4).Value = "'6.75"
$ws.Cells.Item(32, 5).Value = "  -2.64%  "
$ws.Cells.Item(33, 2).Value = "Cosmos"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(33, 4).Value = "'11.37"
$ws.Cells.Item(33, 5).Value = "  -2.56%  "
$ws.Cells.Item(34, 2).Value = "Hedera"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(34, 4).Value = "'0.107"
$ws.Cells.Item(34, 5).Value = "  -3.16%  "
$ws.Cells.Item(35, 2).Value = "OKB"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(35, 4).Value = "'59.74"
$ws.Cells.Item(35, 5).Value = "  -6.29%  "
$ws.Cells.Item(36, 2).Value = "Dai"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(36, 4).Value = "'1.00"
$ws.Cells.Item(36, 5).Value = "  +0.09%  "
$ws.Cells.Item(37, 2).Value = "TheGraph"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(37, 4).Value = "'0.396"
$ws.Cells.Item(37, 5).Value = "  +0.26%  "
$ws.Cells.Item(38, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(38, 4).Value = "'37.28"
$ws.Cells.Item(38, 5).Value = "  -3.42%  "
$ws.Cells.Item(39, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(39, 4).Value = "'0.998"
$ws.Cells.Item(39, 5).Value = "  -0.16%  "
$ws.Cells.Item(40, 4).Value = "0.0₃0714"
$ws.Cells.Item(40, 5).Value = "  +6.86%  "
$ws.Cells.Item(41, 2).Value = "Kaspa"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(41, 4).Value = "'0.126"
$ws.Cells.Item(41, 5).Value = "  -3.33%  "
$ws.Cells.Item(42, 2).Value = "Maker"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(42, 4).Value = "2.941.90"
$ws.Cells.Item(42, 5).Value = "  -3.54%  "
$ws.Cells.Item(43, 2).Value = "Fetch.AI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(43, 4).Value = "'2.52"
$ws.Cells.Item(43, 5).Value = "  +0.01%  "
$ws.Cells.Item(44, 2).Value = "ThetaToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(44, 4).Value = "'2.73"
$ws.Cells.Item(44, 5).Value = "  -7.59%  "
$ws.Cells.Item(45, 2).Value = "VeChain"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(45, 4).Value = "'0.0403"
$ws.Cells.Item(45, 5).Value = "  +1.93%  "
$ws.Cells.Item(46, 2).Value = "WEMIXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(46, 4).Value = "'2.67"
$ws.Cells.Item(46, 5).Value = "  -3.58%  "
$ws.Cells.Item(47, 2).Value = "Stacks"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(47, 4).Value = "'2.82"
$ws.Cells.Item(47, 5).Value = "  +7.92%  "
$ws.Cells.Item(48, 4).Value = "'3.04"
$ws.Cells.Item(48, 5).Value = "  +6.29%  "
$ws.Cells.Item(49, 2).Value = "Stellar"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(49, 4).Value = "'0.127"
$ws.Cells.Item(49, 5).Value = "  +0.68%  "
$ws.Cells.Item(50, 2).Value = "dogwifhat"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(50, 4).Value = "'2.54"
$ws.Cells.Item(50, 5).Value = "  -4.57%  "
$ws.Cells.Item(51, 2).Value = "Monero"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(51, 4).Value = "'135.87"
$ws.Cells.Item(51, 5).Value = "  -1.01%  "
